$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E3: "1234567" -> "123456789"
$ws.Range("E3").Value = "123456789"

# E4: numeric 123456789 -> text "1234567"
$ws.Range("E4").Value = "1234567"

# New row 6: additional test case - copy formatting from row 4 first
$ws.Range("A4:F4").Copy()
$ws.Range("A6:F6").PasteSpecial(-4122)

$ws.Range("A6").Value = "no"
$ws.Range("B6").Value = "test case number 5"
$ws.Range("C6").Value = "Check response on entering blank email and password"
$ws.Range("D6").Value = "dummy"
$ws.Range("E6").Value = "dummy"
$ws.Range("F6").Value = "email or password is invalid"

# Update the active selection to match the target workbook state
$ws.Range("E6").Select()
